# Insert a new data row above current row 72 (shifts rows 72:93 down to 73:94)
# and populate the new row 72 with the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("72:72").Insert()

$ws.Cells.Item(72, 1).Value = 10
$ws.Cells.Item(72, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(72, 3).Value = "La Araucanía"
$ws.Cells.Item(72, 4).Value = 44785
$ws.Cells.Item(72, 5).Value = 9
$ws.Cells.Item(72, 6).Value = "Fruta"
$ws.Cells.Item(72, 7).Value = 100108
$ws.Cells.Item(72, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(72, 9).Value = 100108007
$ws.Cells.Item(72, 10).Value = "Coco"
$ws.Cells.Item(72, 11).Value = "Sin especificar"
$ws.Cells.Item(72, 12).Value = "Primera"
$ws.Cells.Item(72, 13).Value = 15
$ws.Cells.Item(72, 14).Value = 30000
$ws.Cells.Item(72, 15).Value = 30000
$ws.Cells.Item(72, 16).Value = 30000
$ws.Cells.Item(72, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(72, 18).Value = "Perú"
$ws.Cells.Item(72, 19).Value = 1500
$ws.Cells.Item(72, 20).Value = 20
